$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill B:D for rows 1666-1690
$arr = New-Object 'object[,]' 25,3
$arr[0,0] = 71.44
$arr[0,1] = 11
$arr[0,2] = 'whiskybase.com'
$arr[1,0] = 60
$arr[1,1] = 1
$arr[1,2] = 'liquor.com'
$arr[2,0] = 79.94
$arr[2,1] = 81
$arr[2,2] = 'whiskybase.com'
$arr[3,0] = 69
$arr[3,1] = 1
$arr[3,2] = 'whiskyrant.com'
$arr[4,0] = 72
$arr[4,1] = 73
$arr[4,2] = 'whiskybase.com'
$arr[5,0] = 81.75
$arr[5,1] = 4
$arr[5,2] = 'whiskybase.com'
$arr[6,0] = 73
$arr[6,1] = 1
$arr[6,2] = 'whiskyrant.com'
$arr[7,0] = 30
$arr[7,1] = 1
$arr[7,2] = 'breakingbourbon.com'
$arr[8,0] = 75.430000000000007
$arr[8,1] = 64
$arr[8,2] = 'whiskybase.com'
$arr[9,0] = 77.69
$arr[9,1] = 18
$arr[9,2] = 'whiskybase.com'
$arr[10,0] = 67
$arr[10,1] = 1
$arr[10,2] = 'whiskyrant.com'
$arr[11,0] = 64.62
$arr[11,1] = 91
$arr[11,2] = 'whiskybase.com'
$arr[12,0] = 74
$arr[12,1] = 1
$arr[12,2] = 'whiskybase.com'
$arr[13,0] = 63
$arr[13,1] = 2
$arr[13,2] = 'whiskybase.com'
$arr[14,0] = 74
$arr[14,1] = 1
$arr[14,2] = 'balcohols.com'
$arr[15,0] = 73
$arr[15,1] = 1
$arr[15,2] = 'whiskyrant.com'
$arr[16,0] = 76.45
$arr[16,1] = 46
$arr[16,2] = 'whiskybase.com'
$arr[17,0] = 74.37
$arr[17,1] = 45
$arr[17,2] = 'whiskybase.com'
$arr[18,0] = 67.56
$arr[18,1] = 43
$arr[18,2] = 'whiskybase.com'
$arr[19,0] = 73
$arr[19,1] = 1
$arr[19,2] = 'whiskybase.com'
$arr[20,0] = 77
$arr[20,1] = 1
$arr[20,2] = 'whiskybase.com'
$arr[21,0] = 69.37
$arr[21,1] = 93
$arr[21,2] = 'whiskybase.com'
$arr[22,0] = 69.86
$arr[22,1] = 336
$arr[22,2] = 'whiskybase.com'
$arr[23,0] = 74.38
$arr[23,1] = 18
$arr[23,2] = 'whiskybase.com'
$arr[24,0] = 77.2
$arr[24,1] = 66
$arr[24,2] = 'whiskybase.com'
$ws.Range("B1666:D1690").Value = $arr

# Fill B:D for rows 1694-1701
$arr = New-Object 'object[,]' 8,3
$arr[0,0] = 73.260000000000005
$arr[0,1] = 51
$arr[0,2] = 'whiskybase.com'
$arr[1,0] = 72.010000000000005
$arr[1,1] = 278
$arr[1,2] = 'whiskybase.com'
$arr[2,0] = 73.53
$arr[2,1] = 95
$arr[2,2] = 'whiskybase.com'
$arr[3,0] = 67.56
$arr[3,1] = 43
$arr[3,2] = 'whiskybase.com'
$arr[4,0] = 50
$arr[4,1] = 1
$arr[4,2] = 'breakingbourbon.com'
$arr[5,0] = 74.069999999999993
$arr[5,1] = 33
$arr[5,2] = 'whiskybase.com'
$arr[6,0] = 60
$arr[6,1] = 1
$arr[6,2] = 'breakingbourbon.com'
$arr[7,0] = 87
$arr[7,1] = 1
$arr[7,2] = 'whiskybase.com'
$ws.Range("B1694:D1701").Value = $arr

# New rows 1705-1748: A:D
$arr = New-Object 'object[,]' 44,4
$arr[0,0] = 'Black Velvet muovipullo'
$arr[0,1] = 60
$arr[0,2] = 1
$arr[0,3] = 'liquor.com'
$arr[1,0] = 'Cutty Sark Prohibition Edition'
$arr[1,1] = 81.150000000000006
$arr[1,2] = 287
$arr[1,3] = 'whiskybase.com'
$arr[2,0] = 'Jameson muovipullo'
$arr[2,1] = 76
$arr[2,2] = 1
$arr[2,3] = 'whiskyscores.com'
$arr[3,0] = 'The Famous Grouse Smoky Black'
$arr[3,1] = 76.569999999999993
$arr[3,2] = 37
$arr[3,3] = 'whiskybase.com'
$arr[4,0] = 'Four Roses'
$arr[4,1] = 72.72
$arr[4,2] = 294
$arr[4,3] = 'whiskybase.com'
$arr[5,0] = 'Dewar''s Caribbean Smooth Aged 8 Years'
$arr[5,1] = 75.3
$arr[5,2] = 151
$arr[5,3] = 'whiskybase.com'
$arr[6,0] = 'The Famous Grouse Sherry Cask Finish'
$arr[6,1] = 79.88
$arr[6,2] = 29
$arr[6,3] = 'whiskybase.com'
$arr[7,0] = 'The Whistler Triple Oak'
$arr[7,1] = 83.5
$arr[7,2] = 2
$arr[7,3] = 'whiskybase.com'
$arr[8,0] = 'The Busker Triple Cask Triple Smooth'
$arr[8,1] = 78.760000000000005
$arr[8,2] = 44
$arr[8,3] = 'whiskybase.com'
$arr[9,0] = 'Powers Irish Rye'
$arr[9,1] = 83.44
$arr[9,2] = 20
$arr[9,3] = 'whiskybase.com'
$arr[10,0] = 'Spey River Bourbon Cask Single Malt'
$arr[10,1] = 77
$arr[10,2] = 1
$arr[10,3] = 'whiskybase.com'
$arr[11,0] = 'Glen Moray Peated Single Malt'
$arr[11,1] = 79.040000000000006
$arr[11,2] = 294
$arr[11,3] = 'whiskybase.com'
$arr[12,0] = 'Glen Moray Port Cask Finish Single Malt'
$arr[12,1] = 80
$arr[12,2] = 18
$arr[12,3] = 'masterofmalt.com'
$arr[13,0] = 'Glen Moray Sherry Cask Finish Single Malt'
$arr[13,1] = 79.13
$arr[13,2] = 285
$arr[13,3] = 'whiskybase.com'
$arr[14,0] = 'Grant''s Family Reserve'
$arr[14,1] = 66.34
$arr[14,2] = 207
$arr[14,3] = 'whiskybase.com'
$arr[15,0] = 'Tamnavulin Sherry Cask Single Malt'
$arr[15,1] = 80.98
$arr[15,2] = 569
$arr[15,3] = 'whiskybase.com'
$arr[16,0] = 'Islay Mist Mufloni Stout Cask Finish'
$arr[16,1] = 79.75
$arr[16,2] = 6
$arr[16,3] = 'whiskybase.com'
$arr[17,0] = 'Wild Turkey 101'
$arr[17,1] = 86.06
$arr[17,2] = 275
$arr[17,3] = 'whiskybase.com'
$arr[18,0] = 'Jack Daniel''s Old No. 7'
$arr[18,1] = 66.78
$arr[18,2] = 944
$arr[18,3] = 'whiskybase.com'
$arr[19,0] = 'Jack Daniel''s Old No. 7 muovipullo'
$arr[19,1] = 66.78
$arr[19,2] = 944
$arr[19,3] = 'whiskybase.com'
$arr[20,0] = 'Loch Lomond Original Single Malt'
$arr[20,1] = 76.78
$arr[20,2] = 241
$arr[20,3] = 'whiskybase.com'
$arr[21,0] = 'Bowsaw Small Batch Straight Corn American Whiskey'
$arr[21,1] = 78
$arr[21,2] = 17
$arr[21,3] = 'whiskybase.com'
$arr[22,0] = 'Hamiltons Highland Single Malt'
$arr[22,1] = 74.2
$arr[22,2] = 9
$arr[22,3] = 'whiskybase.com'
$arr[23,0] = 'Maker''s Mark'
$arr[23,1] = 50
$arr[23,2] = 1
$arr[23,3] = 'breakingbourbon.com'
$arr[24,0] = 'Tamnavulin White Wine Cask Finish Single Malt'
$arr[24,1] = 81.69
$arr[24,2] = 118
$arr[24,3] = 'whiskybase.com'
$arr[25,0] = 'Jura Red Wine Cask Finish Single Malt'
$arr[25,1] = 80.14
$arr[25,2] = 189
$arr[25,3] = 'whiskybase.com'
$arr[26,0] = 'Tamnavulin German Pinot Noir Cask Finish Single Malt'
$arr[26,1] = 81.349999999999994
$arr[26,2] = 130
$arr[26,3] = 'whiskybase.com'
$arr[27,0] = 'The Singleton 12 Year Old Single Malt'
$arr[27,1] = 79.16
$arr[27,2] = 519
$arr[27,3] = 'whiskybase.com'
$arr[28,0] = 'The Epicurean Blended Malt'
$arr[28,1] = 82
$arr[28,2] = 11
$arr[28,3] = 'whiskybase.com'
$arr[29,0] = 'Glen Moray Elgin Classic Single Malt'
$arr[29,1] = 79.040000000000006
$arr[29,2] = 294
$arr[29,3] = 'whiskybase.com'
$arr[30,0] = 'Ragtime Rye'
$arr[30,1] = 80
$arr[30,2] = 1
$arr[30,3] = 'reddit.com/r/bourbon'
$arr[31,0] = 'Smokey Joe Islay Blended Malt'
$arr[31,1] = 80.510000000000005
$arr[31,2] = 138
$arr[31,3] = 'whiskybase.com'
$arr[32,0] = 'The Quiet Man'
$arr[32,1] = 86
$arr[32,2] = 1
$arr[32,3] = 'reddit.com/r/irishwhiskey'
$arr[33,0] = 'Rittenhouse Straight Rye Whiskey'
$arr[33,1] = 81.819999999999993
$arr[33,2] = 356
$arr[33,3] = 'whiskybase.com'
$arr[34,0] = 'The Dead Rabbit'
$arr[34,1] = 74
$arr[34,2] = 1
$arr[34,3] = 'reddit.com/r/worldwhisky'
$arr[35,0] = 'MacLean''s Nose'
$arr[35,1] = 84
$arr[35,2] = 385
$arr[35,3] = 'whiskybase.com'
$arr[36,0] = 'The Glenlivet Founder''s Reserve Single Malt'
$arr[36,1] = 75.12
$arr[36,2] = 793
$arr[36,3] = 'whiskybase.com'
$arr[37,0] = 'The Legendary Silkie Irish Whiskey'
$arr[37,1] = 71
$arr[37,2] = 1
$arr[37,3] = 'whiskynotes.be'
$arr[38,0] = 'Smoky Scot Islay Single Malt'
$arr[38,1] = 83.34
$arr[38,2] = 164
$arr[38,3] = 'whiskybase.com'
$arr[39,0] = 'Mossburn Speyside Blended Malt'
$arr[39,1] = 81.900000000000006
$arr[39,2] = 143
$arr[39,3] = 'whiskybase.com'
$arr[40,0] = 'Proclamation'
$arr[40,1] = 79.150000000000006
$arr[40,2] = 24
$arr[40,3] = 'whiskybase.com'
$arr[41,0] = 'Johnnie Walker Black Label 12 Year Old'
$arr[41,1] = 75.81
$arr[41,2] = 928
$arr[41,3] = 'whiskybase.com'
$arr[42,0] = 'Naked Malt'
$arr[42,1] = 80
$arr[42,2] = 1
$arr[42,3] = 'reddit.com/r/scotch'
$arr[43,0] = 'Jack Daniel''s Bonded'
$arr[43,1] = 82.3
$arr[43,2] = 82
$arr[43,3] = 'whiskybase.com'
$ws.Range("A1705:D1748").Value = $arr

# New rows 1749-1819: A only
$arr = New-Object 'object[,]' 71,1
$arr[0,0] = 'Powers Three Swallow Single Pot Still'
$arr[1,0] = 'Monkey Shoulder Blended Malt'
$arr[2,0] = 'Chivas Regal 12 Years Old'
$arr[3,0] = 'Glenmorangie X Single Malt'
$arr[4,0] = 'Knob Creek Small Batch Bourbon'
$arr[5,0] = 'BlueWhite Lightly Peated'
$arr[6,0] = 'Legendary Dark Silkie Smoky Irish Whiskey'
$arr[7,0] = 'Redemption Bourbon'
$arr[8,0] = 'Nikka Days'
$arr[9,0] = 'Starward Left-Field Single Malt'
$arr[10,0] = 'Glen Moray Elgin Heritage 12 Year Old Single Malt'
$arr[11,0] = 'The Deacon'
$arr[12,0] = 'Scarabus Batch Strength Single Malt'
$arr[13,0] = 'Hinch Small Batch'
$arr[14,0] = 'Glen Scotia Campbeltown Harbour Single Malt'
$arr[15,0] = 'Laphroaig Select Single Malt'
$arr[16,0] = 'Talisker Storm Single Malt'
$arr[17,0] = 'Glen Garioch 1797 Founder''s Reserve'
$arr[18,0] = 'Glendalough Single Cask Madeira Finish'
$arr[19,0] = 'Glendalough Single Cask Burgundy Finish'
$arr[20,0] = 'Ardbeg Wee Beastie Single Malt'
$arr[21,0] = 'Knob Creek Small Batch Rye'
$arr[22,0] = 'Speyburn 10 Year Old Single Malt'
$arr[23,0] = 'Abasolo Corn Whisky'
$arr[24,0] = 'Compass Box Great King St Glasgow Blend'
$arr[25,0] = 'Teeling Cask Strength Pot Still'
$arr[26,0] = 'Agitator Single Malt Whisky Rök'
$arr[27,0] = 'Peat''s Beast Single Malt'
$arr[28,0] = 'Yellow Rose Premium American'
$arr[29,0] = 'BlueWhite Pääesiintyjät'
$arr[30,0] = 'Woven Superblend'
$arr[31,0] = 'Talisker 10 Year Old Single Malt'
$arr[32,0] = 'Tomintoul 10 Year Old Single Malt'
$arr[33,0] = 'Woodford Reserve Distiller''s Select Rye'
$arr[34,0] = 'Pointe Blanche Single Malt'
$arr[35,0] = 'The English Smokey Single Malt'
$arr[36,0] = 'Hatozaki Pure Malt'
$arr[37,0] = 'Jameson Single Pot Still'
$arr[38,0] = 'Micil Earls Island Single Pot Still'
$arr[39,0] = 'Lindores MCDXCIV Lowland Single Malt'
$arr[40,0] = 'Compass Box Orchard House Blended Malt'
$arr[41,0] = 'Arran Sherry Cask Single Malt'
$arr[42,0] = 'Laphroaig Quarter Cask Single Malt'
$arr[43,0] = 'Dingle Single Malt'
$arr[44,0] = 'Stauning Floor Malted Rye'
$arr[45,0] = 'Spirit of Raskasta Joulua'
$arr[46,0] = 'Yellow Rose Rye'
$arr[47,0] = 'Woodford Reserve Distiller''s Select'
$arr[48,0] = 'Glenmorangie The Original 12 Year Old Single Malt'
$arr[49,0] = 'Peat''s Beast PX Sherry Finish Batch Strength Single Malt'
$arr[50,0] = 'The English Small Batch Rum Cask Single Malt'
$arr[51,0] = 'Johnnie Walker Green Label 15 Year Old Blended Malt'
$arr[52,0] = 'Maker''s Mark Cask Strength'
$arr[53,0] = 'Yellowstone Select'
$arr[54,0] = 'Teeling Blackpitts Peated Single Malt'
$arr[55,0] = 'Lindores Casks of Lindores II Bourbon Lowland Single Malt'
$arr[56,0] = 'Kyrö’s Choice SVP 2025 Wood Smoke Single Vintage 2019'
$arr[57,0] = 'Angel''s Envy'
$arr[58,0] = 'London Square 12 Year Old'
$arr[59,0] = 'Glen Scotia Double Cask Single Malt'
$arr[60,0] = 'Glenrothes 10 Year Old 2014 Signatory Vintage Single Malt'
$arr[61,0] = 'Stauning Smooth & Delicate Double Malt Høst'
$arr[62,0] = 'Jack Daniel''s Old No. 7 muovipullo'
$arr[63,0] = 'Torabhaig Legacy Allt Gleann Single Malt'
$arr[64,0] = 'Tobermory 12 Year Old Single Malt'
$arr[65,0] = 'Smokehead Single Malt'
$arr[66,0] = 'Compass Box The Spice Tree Blended Malt'
$arr[67,0] = 'Kilchoman Machir Bay Single Malt'
$arr[68,0] = 'Kingsbarns Dream to Dram Lowland Single Malt'
$arr[69,0] = 'Penderyn Rich Oak Single Malt'
$arr[70,0] = 'The Whistler P.X. I Love You Single Malt'
$ws.Range("A1749:A1819").Value = $arr

# Column width change for column A
$ws.Columns.Item(1).ColumnWidth = 48.45

# View / selection changes
$ws.Range("B1749").Select()
